$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add the new date entry and description for row 64
$ws.Range("A64").Value = (Get-Date -Year 2013 -Month 4 -Day 8 -Hour 0 -Minute 0 -Second 0 -Millisecond 0)
$ws.Range("B64").Value = "wrote subchapter What is OpenCL? and Components"

# Update the view state: scroll position and active selection
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 40
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B68").Select()
